$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1834319526627219
$ws.Range("C2").Value = 0.5621301775147929
$ws.Range("J2").Value = 0.005917159763313609
$ws.Range("P2").Value = 0.1153846153846154
$ws.Range("S2").Value = 0.1331360946745562
$ws.Range("B3").Value = 0.01015228426395939
$ws.Range("C3").Value = 0.01522842639593909
$ws.Range("J3").Value = 0.02538071065989848
$ws.Range("P3").Value = 0.7868020304568528
$ws.Range("S3").Value = 0.1624365482233502
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.7777777777777778
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.08860759493670886
$ws.Range("D6").Value = 0.008438818565400843
$ws.Range("F6").Value = 0.04641350210970464
$ws.Range("J6").Value = 0.2405063291139241
$ws.Range("O6").Value = 0.02531645569620253
$ws.Range("Q6").Value = 0.1814345991561181
$ws.Range("R6").Value = 0.08860759493670886
$ws.Range("S6").Value = 0.3206751054852321
$ws.Range("B7").Value = 0.09502262443438914
$ws.Range("D7").Value = 0.009049773755656109
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.1583710407239819
$ws.Range("O7").Value = 0.01809954751131222
$ws.Range("Q7").Value = 0.16289592760181
$ws.Range("R7").Value = 0.1040723981900453
$ws.Range("S7").Value = 0.3936651583710407
$ws.Range("B8").Value = 0.1093023255813954
$ws.Range("D8").Value = 0.02093023255813953
$ws.Range("F8").Value = 0.05813953488372093
$ws.Range("J8").Value = 0.08837209302325581
$ws.Range("O8").Value = 0.0186046511627907
$ws.Range("Q8").Value = 0.1930232558139535
$ws.Range("R8").Value = 0.1418604651162791
$ws.Range("S8").Value = 0.3697674418604651
$ws.Range("B9").Value = 0.09405940594059406
$ws.Range("D9").Value = 0.01485148514851485
$ws.Range("F9").Value = 0.07425742574257425
$ws.Range("J9").Value = 0.0891089108910891
$ws.Range("O9").Value = 0.009900990099009901
$ws.Range("Q9").Value = 0.1930693069306931
$ws.Range("R9").Value = 0.1386138613861386
$ws.Range("S9").Value = 0.3861386138613861
$ws.Range("B10").Value = 0.111340206185567
$ws.Range("D10").Value = 0.01993127147766323
$ws.Range("F10").Value = 0.07010309278350516
$ws.Range("J10").Value = 0.1237113402061856
$ws.Range("O10").Value = 0.01168384879725086
$ws.Range("Q10").Value = 0.2268041237113402
$ws.Range("R10").Value = 0.1003436426116839
$ws.Range("S10").Value = 0.3360824742268042
$ws.Range("G11").Value = 0.1457142857142857
$ws.Range("J11").Value = 0.08
$ws.Range("K11").Value = 0.2085714285714286
$ws.Range("L11").Value = 0.5514285714285714
$ws.Range("S11").Value = 0.01428571428571429
$ws.Range("G12").Value = 0.7397959183673469
$ws.Range("J12").Value = 0.2040816326530612
$ws.Range("K12").Value = 0.01530612244897959
$ws.Range("L12").Value = 0.01530612244897959
$ws.Range("S12").Value = 0.02551020408163265
$ws.Range("G13").Value = 0.5600000000000001
$ws.Range("J13").Value = 0.38
$ws.Range("S13").Value = 0.06
$ws.Range("G14").Value = 0.4285714285714285
$ws.Range("J14").Value = 0.1428571428571428
$ws.Range("S14").Value = 0.4285714285714285
$ws.Range("F15").Value = 0.01333333333333333
$ws.Range("H15").Value = 0.1466666666666667
$ws.Range("I15").Value = 0.08888888888888889
$ws.Range("J15").Value = 0.3911111111111111
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.004444444444444444
$ws.Range("N15").Value = 0.004444444444444444
$ws.Range("O15").Value = 0.03555555555555556
$ws.Range("S15").Value = 0.2488888888888889
$ws.Range("F16").Value = 0.0091324200913242
$ws.Range("H16").Value = 0.1415525114155251
$ws.Range("I16").Value = 0.0410958904109589
$ws.Range("J16").Value = 0.4474885844748858
$ws.Range("K16").Value = 0.1187214611872146
$ws.Range("M16").Value = 0.0136986301369863
$ws.Range("N16").Value = 0.0091324200913242
$ws.Range("O16").Value = 0.0547945205479452
$ws.Range("S16").Value = 0.1643835616438356
$ws.Range("F17").Value = 0.007561436672967864
$ws.Range("H17").Value = 0.1436672967863894
$ws.Range("I17").Value = 0.1020793950850662
$ws.Range("J17").Value = 0.4555765595463138
$ws.Range("K17").Value = 0.08317580340264651
$ws.Range("M17").Value = 0.01890359168241966
$ws.Range("N17").Value = 0.001890359168241966
$ws.Range("O17").Value = 0.06994328922495274
$ws.Range("S17").Value = 0.1172022684310019
$ws.Range("F18").Value = 0.01428571428571429
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.06785714285714285
$ws.Range("J18").Value = 0.4392857142857143
$ws.Range("K18").Value = 0.08928571428571429
$ws.Range("M18").Value = 0.02857142857142857
$ws.Range("N18").Value = 0.007142857142857143
$ws.Range("O18").Value = 0.1
$ws.Range("S18").Value = 0.1107142857142857
$ws.Range("F19").Value = 0.0182370820668693
$ws.Range("H19").Value = 0.1899696048632219
$ws.Range("I19").Value = 0.07598784194528875
$ws.Range("J19").Value = 0.3844984802431611
$ws.Range("K19").Value = 0.1231003039513678
$ws.Range("M19").Value = 0.02279635258358663
$ws.Range("N19").Value = 0.002279635258358662
$ws.Range("O19").Value = 0.06155015197568389
$ws.Range("S19").Value = 0.121580547112462
